$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualizei dados bibi: linha 8 (ano 2025) - total_customers, new_customers,
# new_rate e returning_rate recalculados.
$ws.Range("C8").Value = 1171
$ws.Range("E8").Value = 978
$ws.Range("G8").Value = 83.51836037574722
$ws.Range("H8").Value = 16.48163962425278
